$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stash a copy of the "category label" style (currently on A3:A14) into a
# scratch cell far outside the used range, so we still have a reference to it
# after we start overwriting A3:A14 with shifted content.
$ws.Range("A3").Copy()
$ws.Range("Z100").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Row 1: tweak a couple of header-row styles
# ------------------------------------------------------------------
# B1 becomes style s=4 (same definition as G1's style)
$ws.Range("G1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# L1:O1 becomes style s=3 (same as the A1/D1/I1 title style)
$ws.Range("A1").Copy()
$ws.Range("L1:O1").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Shift the "Name/Amount" header + the 12 category rows down by two rows:
#   old A2:B2   (Name/Amount headers)        -> new A4:B4
#   old A3:B14  (Housing..Charity, 12 rows)  -> new A5:B16
# Work from the bottom row upwards so we never clobber a source row before
# it has been copied.
# ------------------------------------------------------------------

# old A3:B14 -> new A5:B16 (12 rows of category/amount)
for ($r = 14; $r -ge 3; $r--) {
    $destRow = $r + 2
    $label = $ws.Range("A$r").Value2
    $ws.Range("A$destRow").Value2 = $label
    $ws.Range("Z100").Copy()
    $ws.Range("A$destRow").PasteSpecial(-4122)
}
$ws.Range("B5:B16").Formula = "=SUMIF(`$G:`$G,A5,`$F:`$F)"

# old A2:B2 (Name/Amount headers) -> new A4:B4
$ws.Range("A4").Value2 = $ws.Range("A2").Value2
$ws.Range("B4").Value2 = $ws.Range("B2").Value2
$ws.Range("D2").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)

# ------------------------------------------------------------------
# New row 2: "Total of sums" / "Total values" labels (no explicit style)
# ------------------------------------------------------------------
$ws.Range("A2").Value = "Total of sums"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "Total values"
$ws.Range("B2").Style = "Normal"

# ------------------------------------------------------------------
# New row 3: roll-up formulas (no explicit style)
# ------------------------------------------------------------------
$ws.Range("A3").Formula = "=SUM(B5:B16)"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Formula = "=SUM(F3:F78)"
$ws.Range("B3").Style = "Normal"

# ------------------------------------------------------------------
# Clean up the scratch cell
# ------------------------------------------------------------------
$ws.Range("Z100").Clear()

# ------------------------------------------------------------------
# Selection: whole new A1:B16 block selected, no distinct active cell
# ------------------------------------------------------------------
$ws.Range("A1:B16").Select()
